$d = $word.ActiveDocument

# Locate the paragraph that contains "LOB1012: Estatística (Requisito
# fraco)" and the paragraph that contains the "© 2020 ..." footer text,
# then delete everything from the end of the former (i.e. right after
# its own paragraph mark) through to the end of the latter (including
# its paragraph mark). This removes the blank paragraph, the
# "Ver no Jupiter..." paragraph and the "© 2020..." paragraph, while
# leaving the paragraph before and the paragraphs after untouched.

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*LOB1012: Estatística (Requisito fraco)*") {
        $startPara = $p
    }
    if ($t -like "*Contact: luizeleno@usp.br*") {
        $endPara = $p
    }
}

$deleteRange = $d.Range($startPara.Range.End, $endPara.Range.End)
$deleteRange.Delete()
